# 02.04.21 cleaning folder structure
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("analysis")

# Update the timestamp label in A1
$ws.Range("A1").Value = "02/04/2021 11:46"

# Update C3/C4 (totals)
$ws.Range("C3").Value = 0.4326052928353927
$ws.Range("C4").Value = 0.5673947071646074

# Update C column (series "c")
$ws.Range("C14").Value = 0
$ws.Range("C15").Value = 0.6499526708767375
$ws.Range("C16").Value = 0.679930242599976
$ws.Range("C17").Value = 0.7347938574539591
$ws.Range("C18").Value = 0.5177154407963372
$ws.Range("C20").Value = 0.6854190159458917
$ws.Range("C21").Value = 0.7877297936304034
$ws.Range("C22").Value = 0.8523740931847611
$ws.Range("C23").Value = 0.9436683563245826
$ws.Range("C24").Value = 0.895664937390644
$ws.Range("C25").Value = 1
$ws.Range("C26").Value = 1
$ws.Range("C27").Value = 1

# Update D column (series "y")
$ws.Range("D16").Value = 0.227608881117836
$ws.Range("D17").Value = 0.2427359598426374
$ws.Range("D18").Value = 0.3079322317258981
$ws.Range("D19").Value = 0.2489940375985687
$ws.Range("D20").Value = 0.2181330311580482
$ws.Range("D21").Value = 0.1894773278565225
$ws.Range("D22").Value = 0.106548904989473
$ws.Range("D23").Value = 0.01490999286018015
$ws.Range("D24").Value = 0
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 0
$ws.Range("D27").Value = 0

$wb.Save()
